$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = "[Python - 프로그래머스] 두 원 사이의 정수 쌍"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-%EB%91%90-%EC%9B%90-%EC%82%AC%EC%9D%B4%EC%9D%98-%EC%A0%95%EC%88%98-%EC%8C%8D"

# Row 23
$ws.Range("D23").Value = "[공개][np.stack설명] 넘파이numpy의  stack에 대한 graphical 설명"
$ws.Range("E23").Value = "https://theonly1.tistory.com/3177"

# Row 28
$ws.Range("D28").Value = "ROS:: rosdep, wstool 명령어의 역할과 package 의존성"
$ws.Range("E28").Value = "https://ropiens.tistory.com/222"

# Row 30
$ws.Range("D30").Value = "John-analyst"

# Row 51
$ws.Range("D51").Value = "[python] 추상 클래스(abstract class) 이해하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EC%B6%94%EC%83%81-%ED%81%B4%EB%9E%98%EC%8A%A4abstract-class-%EC%9D%B4%ED%95%B4%ED%95%98%EA%B8%B0"
